$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '29.083.10'
Set-TextValue 'E2' '  +1.53%  '
Set-TextValue 'D3' '1.930.38'
Set-TextValue 'E3' '  +2.36%  '
Set-TextValue 'D4' '1.004'
Set-TextValue 'E4' '  -0.04%  '
Set-TextValue 'D5' '325.81'
Set-TextValue 'E5' '  +1.07%  '
Set-TextValue 'D6' '1.002'
Set-TextValue 'E6' '  -0.05%  '
Set-TextValue 'D7' '0.4607'
Set-TextValue 'E7' '  +0.93%  '
Set-TextValue 'E8' '  +0.99%  '
Set-TextValue 'D9' '0.07744'
Set-TextValue 'E9' '  +0.56%  '
Set-TextValue 'D10' '0.9828'
Set-TextValue 'E10' '  +2.31%  '
Set-TextValue 'D11' '22.53'
Set-TextValue 'E11' '  +2.72%  '
Set-TextValue 'D12' '1.940.55'
Set-TextValue 'E12' '  +4.26%  '
Set-TextValue 'D13' '6.980'
Set-TextValue 'E13' '  +0.47%  '
Set-TextValue 'E14' '  +0.74%  '
Set-TextValue 'D15' '0.07020'
Set-TextValue 'E15' '  -0.33%  '
Set-TextValue 'B16' 'Litecoin'
Set-TextValue 'C16' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D16' '84.49'
Set-TextValue 'E16' '  +1.78%  '
Set-TextValue 'B17' 'BinanceUSD'
Set-TextValue 'C17' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D17' '1.003'
Set-TextValue 'E17' '  -0.14%  '
Set-TextValue 'D18' '0.000009539'
Set-TextValue 'E18' '  +0.48%  '
Set-TextValue 'D19' '16.78'
Set-TextValue 'E19' '  +0.52%  '
Set-TextValue 'D20' '1.002'
Set-TextValue 'E20' '  -0.11%  '
Set-TextValue 'D21' '29.117.38'
Set-TextValue 'E21' '  +1.70%  '
Set-TextValue 'D22' '5.353'
Set-TextValue 'E22' '  +0.22%  '
Set-TextValue 'D23' '10.98'
Set-TextValue 'E23' '  +1.25%  '
Set-TextValue 'D24' '2.081'
Set-TextValue 'E24' '  +1.00%  '
Set-TextValue 'D25' '157.91'
Set-TextValue 'E25' '  +1.75%  '
Set-TextValue 'D26' '19.07'
Set-TextValue 'E26' '  +0.26%  '
Set-TextValue 'D27' '5.699'
Set-TextValue 'E27' '  +1.35%  '
Set-TextValue 'D28' '118.02'
Set-TextValue 'E28' '  +0.91%  '
Set-TextValue 'D29' '1.855'
Set-TextValue 'E29' '  +1.92%  '
Set-TextValue 'D30' '0.09342'
Set-TextValue 'E30' '  +1.20%  '
Set-TextValue 'D31' '0.8677'
Set-TextValue 'E31' '  +2.23%  '
Set-TextValue 'E32' '  +1.30%  '
Set-TextValue 'D33' '1.250'
Set-TextValue 'E33' '  +0.42%  '
Set-TextValue 'D34' '3.028'
Set-TextValue 'E34' '  -0.82%  '
Set-TextValue 'D35' '0.05718'
Set-TextValue 'E35' '  +1.62%  '
Set-TextValue 'E36' '  +0.92%  '
Set-TextValue 'D37' '1.002'
Set-TextValue 'E37' '  -0.03%  '
Set-TextValue 'E38' '  +0.89%  '
Set-TextValue 'D39' '3.066'
Set-TextValue 'E39' '  +14.09%  '
Set-TextValue 'D40' '7.561'
Set-TextValue 'E40' '  +1.77%  '
Set-TextValue 'D41' '0.5533'
Set-TextValue 'E41' '  +1.10%  '
Set-TextValue 'E42' '  +0.48%  '
Set-TextValue 'D43' '0.000002913'
Set-TextValue 'E43' '  -0.48%  '
Set-TextValue 'D44' '9.391'
Set-TextValue 'E44' '  +2.07%  '
Set-TextValue 'D45' '2.237'
Set-TextValue 'E45' '  +7.56%  '
Set-TextValue 'D46' '0.5208'
Set-TextValue 'E46' '  +0.95%  '
Set-TextValue 'E47' '  -0.15%  '
Set-TextValue 'D48' '0.06927'
Set-TextValue 'E48' '  +2.55%  '
Set-TextValue 'D49' '1.784'
Set-TextValue 'E49' '  +0.91%  '
Set-TextValue 'D50' '110.58'
Set-TextValue 'E50' '  +0.03%  '
Set-TextValue 'D51' '70.22'
Set-TextValue 'E51' '  +1.21%  '
